$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "67.077.00"
$ws.Range("E2").Value = "  -0.70%  "

$ws.Range("D3").Value = "2.613.01"
$ws.Range("E3").Value = "  -1.23%  "

Set-TextValue "D5" "588.38"
$ws.Range("E5").Value = "  -1.66%  "

Set-TextValue "D6" "164.86"
$ws.Range("E6").Value = "  -2.04%  "

Set-TextValue "D8" "0.530"
$ws.Range("E8").Value = "  -2.51%  "

$ws.Range("D9").Value = "2.612.26"
$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("E10").Value = "  -5.14%  "

$ws.Range("E11").Value = "  +0.81%  "

Set-TextValue "D12" "0.364"
$ws.Range("E12").Value = "  -0.55%  "

Set-TextValue "D13" "5.21"
$ws.Range("E13").Value = "  -0.63%  "

Set-TextValue "D14" "27.23"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("D15").Value = "3.089.73"
$ws.Range("E15").Value = "  -0.93%  "

$ws.Range("E16").Value = "  -2.95%  "

$ws.Range("D17").Value = "66.992.98"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "2.619.43"
$ws.Range("E18").Value = "  -0.84%  "

Set-TextValue "D19" "11.73"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("E20").Value = "  -1.97%  "

Set-TextValue "D21" "354.99"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("E22").Value = "  -3.14%  "

$ws.Range("E23").Value = "  -3.57%  "

Set-TextValue "D24" "10.46"
$ws.Range("E24").Value = "  -4.81%  "

Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("E26").Value = "  -4.98%  "

Set-TextValue "D27" "69.27"
$ws.Range("E27").Value = "  -1.92%  "

$ws.Range("D28").Value = "2.747.15"
$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("D30").Value = "0.0₃0993"
$ws.Range("E30").Value = "  -3.39%  "

Set-TextValue "D31" "542.95"
$ws.Range("E31").Value = "  -2.50%  "

$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("E33").Value = "  -4.27%  "

$ws.Range("E34").Value = "  -3.16%  "

$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("E37").Value = "  -4.69%  "

Set-TextValue "D38" "157.39"
$ws.Range("E38").Value = "  +0.21%  "

$ws.Range("E39").Value = "  -2.79%  "

$ws.Range("E40").Value = "  -2.61%  "

$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("E42").Value = "  -2.13%  "

Set-TextValue "D43" "5.12"
$ws.Range("E43").Value = "  -3.34%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  -5.12%  "

$ws.Range("E46").Value = "  -1.91%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D47" "0.576"
$ws.Range("E47").Value = "  -3.58%  "

$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D48" "150.87"
$ws.Range("E48").Value = "  -1.87%  "

$ws.Range("E49").Value = "  -2.89%  "

$ws.Range("E50").Value = "  -1.96%  "

Set-TextValue "D51" "0.0768"
$ws.Range("E51").Value = "  -1.51%  "

